$d = $word.ActiveDocument

# Locate the paragraph that ends with "Mi segunda línea de código " and
# insert a brand-new paragraph right after it (before the trailing blank
# paragraph), carrying the same "es-ES" language formatting.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Mi segunda línea de código*") {
        $target = $p
        break
    }
}

$target.Range.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Range.Text = "Nueva línea de código de la rama master "
